$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.706.46'
$ws.Range("E2").Value = '  -0.62%  '

$ws.Range("D3").Value = '2.577.65'
$ws.Range("E3").Value = '  +1.04%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = '''581.43'
$ws.Range("E5").Value = '  -0.76%  '

$ws.Range("D6").Value = '''145.09'
$ws.Range("E6").Value = '  -1.60%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("E8").Value = '  +1.24%  '

$ws.Range("E9").Value = '  +0.87%  '

$ws.Range("E10").Value = '  +0.88%  '

$ws.Range("D11").Value = '''0.151'
$ws.Range("E11").Value = '  -0.19%  '

$ws.Range("E12").Value = '  -0.90%  '

$ws.Range("E13").Value = '  -1.68%  '

$ws.Range("D14").Value = '3.039.22'
$ws.Range("E14").Value = '  +1.14%  '

$ws.Range("D15").Value = '62.592.65'
$ws.Range("E15").Value = '  -0.67%  '

$ws.Range("E16").Value = '  +1.04%  '

$ws.Range("D17").Value = '2.580.61'
$ws.Range("E17").Value = '  +0.84%  '

$ws.Range("E18").Value = '  -1.26%  '

$ws.Range("D19").Value = '''339.97'
$ws.Range("E19").Value = '  +0.81%  '

$ws.Range("E20").Value = '  +0.67%  '

$ws.Range("D21").Value = '''6.66'
$ws.Range("E21").Value = '  -1.80%  '

$ws.Range("E22").Value = '  -0.03%  '

$ws.Range("D23").Value = '''67.34'
$ws.Range("E23").Value = '  +2.39%  '

$ws.Range("D24").Value = '2.706.23'
$ws.Range("E24").Value = '  +0.99%  '

$ws.Range("D25").Value = '''0.165'
$ws.Range("E25").Value = '  -2.54%  '

$ws.Range("E26").Value = '  -2.56%  '

$ws.Range("D27").Value = '''1.00'
$ws.Range("E27").Value = '  +0.34%  '

$ws.Range("D28").Value = '''7.87'
$ws.Range("E28").Value = '  +1.44%  '

$ws.Range("E29").Value = '  -1.96%  '

$ws.Range("E30").Value = '  -1.55%  '

$ws.Range("E31").Value = '  -4.10%  '

$ws.Range("D32").Value = '0.0₃0815'
$ws.Range("E32").Value = '  -0.28%  '

$ws.Range("D33").Value = '''460.62'
$ws.Range("E33").Value = '  +9.65%  '

$ws.Range("D34").Value = '''176.47'
$ws.Range("E34").Value = '  -1.20%  '

$ws.Range("D35").Value = '''1.59'
$ws.Range("E35").Value = '  +2.69%  '

$ws.Range("E36").Value = '  -0.01%  '

$ws.Range("E37").Value = '  -1.05%  '

$ws.Range("E38").Value = '  -1.17%  '

$ws.Range("E39").Value = '  +2.89%  '

$ws.Range("E40").Value = '  +0.00%  '

$ws.Range("E41").Value = '  -3.32%  '

$ws.Range("D42").Value = '''157.87'
$ws.Range("E42").Value = '  +4.98%  '

$ws.Range("E43").Value = '  -1.62%  '

$ws.Range("D44").Value = '''21.05'
$ws.Range("E44").Value = '  +0.81%  '

$ws.Range("E45").Value = '  +4.02%  '

$ws.Range("E46").Value = '  -1.10%  '

$ws.Range("E47").Value = '  -0.98%  '

$ws.Range("E48").Value = '  -1.34%  '

$ws.Range("D49").Value = '''18.16'
$ws.Range("E49").Value = '  -0.95%  '

$ws.Range("D50").Value = '''11.41'
$ws.Range("E50").Value = '  +0.79%  '

$ws.Range("D51").Value = '''1.71'
$ws.Range("E51").Value = '  -0.72%  '
